# Implement basket-based elective scheduling with common time slots across all branches
$wb = $excel.ActiveWorkbook

$wsA = $wb.Worksheets.Item("Section_A")
$wsB = $wb.Worksheets.Item("Section_B")

# --- Section_A updates ---
$wsA.Range("B2").Value = "CS161"
$wsA.Range("C2").Value = "HS161"
$wsA.Range("E2").Value = "MA162"

$wsA.Range("B3").Value = "Free"
$wsA.Range("C3").Value = "CS161"
$wsA.Range("D3").Value = "EC161"
$wsA.Range("E3").Value = "Free"

$wsA.Range("B5").Value = "EC161"
$wsA.Range("C5").Value = "MA162"
$wsA.Range("D5").Value = "MA161"
$wsA.Range("E5").Value = "HS161"
$wsA.Range("F5").Value = "DS161"

$wsA.Range("B7").Value = "MA161"
$wsA.Range("C7").Value = "DS161"
$wsA.Range("D7").Value = "CS161"
$wsA.Range("E7").Value = "EC161"
$wsA.Range("F7").Value = "EC161"

# --- Section_B updates ---
$wsB.Range("B2").Value = "EC161"
$wsB.Range("C2").Value = "HS161"
$wsB.Range("D2").Value = "Free"
$wsB.Range("E2").Value = "HS161"

$wsB.Range("B3").Value = "MA161"
$wsB.Range("C3").Value = "CS161"
$wsB.Range("D3").Value = "CS161"
$wsB.Range("E3").Value = "MA161"

$wsB.Range("B5").Value = "CS161"
$wsB.Range("C5").Value = "DS161"
$wsB.Range("D5").Value = "HS161"
$wsB.Range("E5").Value = "MA162"
$wsB.Range("F5").Value = "MA162"

$wsB.Range("B7").Value = "Free"
$wsB.Range("C7").Value = "EC161"
$wsB.Range("D7").Value = "EC161"
$wsB.Range("E7").Value = "DS161"
$wsB.Range("F7").Value = "DS161"
